# Add season record columns (Wins / Losses / Ties) to the right of the
# existing "Unnamed: 28" column (AC), for the Sheet1 roster table that
# spans rows 1-50 (A1:AC50 -> A1:AF50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# formatting used by the rest of row 1, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-50) gets the same season record values.
$ws.Range("AD2:AD50").Value = 79
$ws.Range("AE2:AE50").Value = 83
$ws.Range("AF2:AF50").Value = 0
